$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 39
$ws.Range("H39").Value = 181.83333
$ws.Range("I39").Value = 181.83333
$ws.Range("K39").Value = 545.49999
$ws.Range("M39").Value = -249.49999
# Row 51
$ws.Range("H51").Value = 12935.143
$ws.Range("I51").Value = 7500
$ws.Range("K51").Value = 7500
$ws.Range("M51").Value = -7016
# Row 69
$ws.Range("H69").Value = 21122.309
$ws.Range("I69").Value = 21508.182
$ws.Range("J69").Value = 19000
$ws.Range("K69").Value = 64524.546
$ws.Range("L69").Value = 57000
$ws.Range("M69").Value = -63650.546
$ws.Range("N69").Value = -58748
# Row 72
$ws.Range("H72").Value = 21122.309
$ws.Range("I72").Value = 21508.182
$ws.Range("J72").Value = 19000
$ws.Range("K72").Value = 193573.638
$ws.Range("L72").Value = 171000
$ws.Range("M72").Value = -189205.638
$ws.Range("N72").Value = -179736
# Row 132
$ws.Range("H132").Value = 237419.56
$ws.Range("I132").Value = 316117.66
$ws.Range("J132").Value = 6946.5713
$ws.Range("K132").Value = 948352.98
$ws.Range("L132").Value = 20839.7139
$ws.Range("M132").Value = -945822.98
$ws.Range("N132").Value = -25899.7139
# Row 135
$ws.Range("H135").Value = 4845.6216
$ws.Range("I135").Value = 1875.5172
$ws.Range("K135").Value = 16879.6548
$ws.Range("M135").Value = -14344.6548
# Row 137
$ws.Range("H137").Value = 10537.36
$ws.Range("I137").Value = 6482.0586
$ws.Range("K137").Value = 19446.1758
$ws.Range("M137").Value = -16896.1758
# Row 138
$ws.Range("H138").Value = 4072.7708
$ws.Range("I138").Value = 463.30768
$ws.Range("J138").Value = 5413.4287
$ws.Range("K138").Value = 1389.92304
$ws.Range("L138").Value = 16240.2861
$ws.Range("M138").Value = 3750.07696
$ws.Range("N138").Value = -26520.2861

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 330.18182
$ws.Range("I4").Value = 280.375
$ws.Range("J4").Value = 463
$ws.Range("K4").Value = 280.375
$ws.Range("L4").Value = 463
$ws.Range("M4").Value = -164.375
$ws.Range("N4").Value = -695
# Row 5
$ws.Range("H5").Value = 148
$ws.Range("I5").Value = 148
$ws.Range("K5").Value = 148
$ws.Range("M5").Value = -36
# Row 17
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()
# Row 32
$ws.Range("H32").Value = 1341009.1
$ws.Range("I32").Value = 1126.2
$ws.Range("J32").Value = 8555764
$ws.Range("K32").Value = 1126.2
$ws.Range("L32").Value = 8555764
$ws.Range("M32").Value = -839.2
$ws.Range("N32").Value = -8556338
# Row 61
$ws.Range("H61").Value = 4255.552
$ws.Range("I61").Value = 4312.909
$ws.Range("J61").Value = 4075.2856
$ws.Range("K61").Value = 4312.909
$ws.Range("L61").Value = 4075.2856
$ws.Range("M61").Value = -4100.909
$ws.Range("N61").Value = -4499.2856
# Row 74
$ws.Range("H74").Value = 4293.35
$ws.Range("I74").Value = 5207.2856
$ws.Range("J74").Value = 2160.8333
$ws.Range("K74").Value = 5207.2856
$ws.Range("L74").Value = 2160.8333
$ws.Range("M74").Value = -4333.2856
$ws.Range("N74").Value = -3908.8333
# Row 77
$ws.Range("H77").Value = 4293.35
$ws.Range("I77").Value = 5207.2856
$ws.Range("J77").Value = 2160.8333
$ws.Range("K77").Value = 26036.428
$ws.Range("L77").Value = 10804.1665
$ws.Range("M77").Value = -21668.428
$ws.Range("N77").Value = -19540.1665
# Row 132
$ws.Range("H132").Value = 668222.8
$ws.Range("I132").Value = 843357.75
$ws.Range("J132").Value = 95053.91
$ws.Range("K132").Value = 2530073.25
$ws.Range("L132").Value = 285161.73
$ws.Range("M132").Value = -2527543.25
$ws.Range("N132").Value = -290221.73
# Row 136
$ws.Range("H136").Value = 4255.552
$ws.Range("I136").Value = 4312.909
$ws.Range("J136").Value = 4075.2856
$ws.Range("K136").Value = 12938.727
$ws.Range("L136").Value = 12225.8568
$ws.Range("M136").Value = -10388.727
$ws.Range("N136").Value = -17325.8568

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 148
$ws.Range("I4").Value = 148
$ws.Range("K4").Value = 148
$ws.Range("M4").Value = -33
# Row 22
$ws.Range("H22").Value = 5555.049
$ws.Range("I22").Value = 2576.5557
$ws.Range("K22").Value = 2576.5557
$ws.Range("M22").Value = -2403.5557
# Row 134
$ws.Range("H134").Value = 1254962.6
$ws.Range("I134").Value = 1474176.6
$ws.Range("J134").Value = 12749.833
$ws.Range("K134").Value = 4422529.800000001
$ws.Range("L134").Value = 38249.499
$ws.Range("M134").Value = -4419994.800000001
$ws.Range("N134").Value = -43319.499

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4130.4287
$ws.Range("I31").Value = 1416.909
$ws.Range("K31").Value = 1416.909
$ws.Range("M31").Value = -1121.909
# Row 34
$ws.Range("H34").Value = 4130.4287
$ws.Range("I34").Value = 1416.909
$ws.Range("K34").Value = 1416.909
$ws.Range("M34").Value = -1214.909
# Row 119
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
# Row 132
$ws.Range("H132").Value = 4159.9375
$ws.Range("I132").Value = 3427.7837
$ws.Range("J132").Value = 6622.636
$ws.Range("K132").Value = 10283.3511
$ws.Range("L132").Value = 19867.908
$ws.Range("M132").Value = -7753.3511
$ws.Range("N132").Value = -24927.908

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 75
$ws.Range("H75").Value = 1813.3334
$ws.Range("J75").Value = 2495
$ws.Range("L75").Value = 7485
$ws.Range("N75").Value = -9481
# Row 78
$ws.Range("H78").Value = 1813.3334
$ws.Range("J78").Value = 2495
$ws.Range("L78").Value = 22455
$ws.Range("N78").Value = -32439
# Row 114
$ws.Range("H114").Value = 68175.39999999999
$ws.Range("I114").Value = 1705.625
$ws.Range("J114").Value = 144140.86
$ws.Range("K114").Value = 5116.875
$ws.Range("L114").Value = 432422.58
$ws.Range("M114").Value = -1862.875
$ws.Range("N114").Value = -438930.58
# Row 128
$ws.Range("H128").Value = 154333
$ws.Range("I128").Value = 154333
$ws.Range("K128").Value = 462999
$ws.Range("M128").Value = -458019
# Row 137
$ws.Range("H137").Value = 1939
$ws.Range("I137").Value = 1152.3334
$ws.Range("J137").Value = 2368.0908
$ws.Range("K137").Value = 3457.0002
$ws.Range("L137").Value = 7104.2724
$ws.Range("M137").Value = 1642.9998
$ws.Range("N137").Value = -17304.2724

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 7
$ws.Range("H7").Value = 17999.8
$ws.Range("I7").Value = 19999.666
$ws.Range("J7").Value = 15000
$ws.Range("K7").Value = 19999.666
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = -19887.666
$ws.Range("N7").Value = -15224
# Row 8
$ws.Range("H8").Value = 17999.8
$ws.Range("I8").Value = 19999.666
$ws.Range("J8").Value = 15000
$ws.Range("K8").Value = 19999.666
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = -19860.666
$ws.Range("N8").Value = -15278
# Row 80
$ws.Range("H80").Value = 5765.8667
$ws.Range("I80").Value = 4419.7144
$ws.Range("J80").Value = 6943.75
$ws.Range("K80").Value = 4419.7144
$ws.Range("L80").Value = 6943.75
$ws.Range("M80").Value = -3421.7144
$ws.Range("N80").Value = -8939.75
# Row 83
$ws.Range("H83").Value = 5765.8667
$ws.Range("I83").Value = 4419.7144
$ws.Range("J83").Value = 6943.75
$ws.Range("K83").Value = 22098.572
$ws.Range("L83").Value = 34718.75
$ws.Range("M83").Value = -17106.572
$ws.Range("N83").Value = -44702.75
# Row 107
$ws.Range("H107").Value = 604.75
$ws.Range("J107").Value = 663.7143
$ws.Range("L107").Value = 663.7143
$ws.Range("N107").Value = -4503.7143
# Row 122
$ws.Range("H122").Value = 3936.6956
$ws.Range("I122").Value = 2939.5334
$ws.Range("J122").Value = 5806.375
$ws.Range("K122").Value = 8818.600199999999
$ws.Range("L122").Value = 17419.125
$ws.Range("M122").Value = -6368.600199999999
$ws.Range("N122").Value = -22319.125
# Row 126
$ws.Range("H126").Value = 16133748
$ws.Range("I126").Value = 25003100
$ws.Range("K126").Value = 75009300
$ws.Range("M126").Value = -75006830
# Row 132
$ws.Range("H132").Value = 15876351
$ws.Range("I132").Value = 23812644
$ws.Range("J132").Value = 3765.1428
$ws.Range("K132").Value = 71437932
$ws.Range("L132").Value = 11295.4284
$ws.Range("M132").Value = -71435402
$ws.Range("N132").Value = -16355.4284

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
# Row 3
$ws.Range("H3").Value = 2500
$ws.Range("I3").Value = 2500
$ws.Range("K3").Value = 2500
$ws.Range("M3").Value = -2388
# Row 15
$ws.Range("H15").Value = 2500
$ws.Range("I15").Value = 2500
$ws.Range("K15").Value = 2500
$ws.Range("M15").Value = -2330
# Row 55
$ws.Range("H55").Value = 2763.1333
$ws.Range("I55").Value = 1000.82355
$ws.Range("J55").Value = 5067.6924
$ws.Range("K55").Value = 1000.82355
$ws.Range("L55").Value = 5067.6924
$ws.Range("M55").Value = -827.82355
$ws.Range("N55").Value = -5413.6924
# Row 132
$ws.Range("H132").Value = 6162.6943
$ws.Range("I132").Value = 5370.8076
$ws.Range("K132").Value = 16112.4228
$ws.Range("M132").Value = -13582.4228

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 7401.5
$ws.Range("I132").Value = 6584.2354
$ws.Range("K132").Value = 19752.7062
$ws.Range("M132").Value = -17222.7062
